# "role change" - the id column's declared type (row 2) is switched
# from "String" to "int".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "int"

# Final selection lands on B8 (matches the saved view state in the diff).
$ws.Range("B8").Select()
